$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 11

# --- A11: year label, styled like the other year cells in column A (bold, centered, boxed) ---
# Copy the style from A10 (the previous year row) so the same cellXf is reused
# rather than a brand-new style being created.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item($row, 1).Value = "2021年"

# --- B11..AQ11: numeric data for 2021 ---
$values = @{
    2  = 4671.88
    3  = 1079.8
    4  = 182.43
    6  = 3538.74
    7  = 8248.74
    8  = 1094.98
    9  = 4082.8
    10 = 556.12
    11 = 148646.06
    12 = 776.24
    13 = 129.88
    14 = 61.66
    15 = 1412.46
    16 = 3045.86
    17 = 179.02
    18 = 222.98
    19 = 3027.17
    20 = 1528.72
    21 = 22699.77
    23 = 1975.96
    24 = 2635.19
    25 = 7610.41
    26 = 9320.610000000001
    27 = 1097.9
    28 = 2111.17
    29 = 2314.92
    30 = 1878.07
    31 = 1754.57
    32 = 33806.85
    33 = 7250.9
    34 = 3341.69
    35 = 1541.17
    36 = 284.74
    37 = 3959.79
    38 = 1688.21
    39 = 3203.69
    40 = 63.35
    41 = 2993.55
    42 = 3146.17
    43 = 121.73
}

foreach ($col in $values.Keys) {
    $ws.Cells.Item($row, $col).Value = $values[$col]
}

# --- E11 and V11: present-but-empty cells (those columns have no data for this row,
#     matching the empty placeholder cells used throughout the sheet, e.g. E9/E10/V9/V10) ---
$blankStyle = $ws.Cells.Item(1000, 500).Style
foreach ($col in 5, 22) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'"
    $c.Style = $blankStyle
}
